# Insert a new data row at row 349 (pushing rows 349:401 down to 350:402)
# and populate it with a new "Cebollín" price observation for
# "Vega Modelo de Temuco" (Región Metropolitana, 2022-08-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 349..401 down to 350..402, creating a blank row 349.
$ws.Rows.Item(349).Insert()

# Fill in the new row 349 with the new record.
$ws.Cells.Item(349, 1).Value = 10
$ws.Cells.Item(349, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(349, 3).Value = "La Araucanía"
$ws.Cells.Item(349, 4).Value = 44776
$ws.Cells.Item(349, 5).Value = 9
$ws.Cells.Item(349, 6).Value = 100112037
$ws.Cells.Item(349, 7).Value = "Cebollín"
$ws.Cells.Item(349, 8).Value = "Sin especificar"
$ws.Cells.Item(349, 9).Value = "Primera"
$ws.Cells.Item(349, 10).Value = 50
$ws.Cells.Item(349, 11).Value = 8000
$ws.Cells.Item(349, 12).Value = 8000
$ws.Cells.Item(349, 13).Value = 8000
$ws.Cells.Item(349, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(349, 15).Value = "Región Metropolitana"
$ws.Cells.Item(349, 16).Value = 667
$ws.Cells.Item(349, 17).Value = 12
$ws.Cells.Item(349, 18).Value = "Hortaliza"
